# "Generate Report for Handback"
#
# The localization status report is regenerated: the zh-cn and de-de files
# have been handed back and are now in sync with en-US, so their Status
# and Latest Handback DateTime are refreshed and the stale "handback file
# is not latest" error is cleared. Column widths are re-fit to the new
# (longer Status / emptied Error Detail) cell content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# --- zh-cn detail sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-27 22:49:24"
$zhcn.Range("P2").Value = ""

# --- de-de detail sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-27 22:49:30"
$dede.Range("P2").Value = ""

# --- Re-fit column widths to match the new cell contents ---
# Status column widened (longer text), Error Detail column narrowed
# (now empty) back down toward its header width.
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(16).ColumnWidth = 13.0

$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(16).ColumnWidth = 13.0
